$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price / volume(1h) values as described in the commit diff.
# Numeric-looking Price values are prefixed with a literal leading apostrophe and
# restyled to "Normal" afterwards so Excel stores them as plain text (matching the
# original inlineStr cell type) instead of auto-coercing them into floating point
# numbers, which would lose formatting such as trailing zeros.

$ws.Range("D2").Value = "62.836.34"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "2.464.97"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'571.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "'147.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("D11").Value = "'5.18"
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").Value = "'28.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").Value = "'0.0000174"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "

$ws.Range("D16").Value = "62.608.70"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "2.466.94"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "'7.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.03%  "

$ws.Range("E19").Value = "  -2.89%  "

$ws.Range("D20").Value = "'2.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").Value = "'322.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("E24").Value = "  +3.47%  "

$ws.Range("D25").Value = "'64.90"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'644.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.12%  "

$ws.Range("D28").Value = "0.0₃0964"
$ws.Range("E28").Value = "  -3.44%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("E30").Value = "  -3.56%  "

$ws.Range("D31").Value = "'7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  -3.77%  "

$ws.Range("E36").Value = "  -2.92%  "

$ws.Range("D37").Value = "'5.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "

$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("D39").Value = "'150.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("D41").Value = "'2.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("E42").Value = "  -2.12%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "0.0₆0304"
$ws.Range("E44").Value = "  -3.31%  "

$ws.Range("D45").Value = "'153.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.88%  "

$ws.Range("D46").Value = "'15.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").Value = "'20.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").Value = "'0.0903"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
